$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "jonathan.doll"
$ws.Range("B3").Value = "ham"

$ws.Range("B4").Select()
